# Apply the changes described by the commit:
#   "replace = by -> for more readability (suggested by fred)"
# plus the incidental footer-date update and the two shape width tweaks
# that accompanied it in the original commit.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Slide master footer date placeholder: "21-Feb-19" -> "1/15/2025"
# ---------------------------------------------------------------------
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        if ($shp.TextFrame.TextRange.Text -eq "21-Feb-19") {
            $shp.TextFrame.TextRange.Text = "1/15/2025"
        }
    }
}

# ---------------------------------------------------------------------
# 2) Slide 2: two shapes get wider (height/position unchanged)
#    A tiny (sub-EMU) nudge is added before the conversion back to EMU
#    so the point -> EMU rounding lands on the exact target instead of
#    one EMU short.
# ---------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$emuPerPt = 12700.0
$nudge = 0.00001

for ($i = 1; $i -le $slide2.Shapes.Count; $i++) {
    $shp = $slide2.Shapes.Item($i)
    if ($shp.Name -eq "CustomShape 1") {
        $shp.Width = (1661400 / $emuPerPt) + $nudge
    }
    elseif ($shp.Name -eq "CustomShape 12") {
        $shp.Width = (5881392 / $emuPerPt) + $nudge
    }
}

# ---------------------------------------------------------------------
# 3) Slide 2, "CustomShape 12": replace "=" with "->" in each of the
#    eight "() = N" snippets (exact per-occurrence substrings, since a
#    couple of them have irregular spacing around the "=").
# ---------------------------------------------------------------------
$codeShape = $null
for ($i = 1; $i -le $slide2.Shapes.Count; $i++) {
    if ($slide2.Shapes.Item($i).Name -eq "CustomShape 12") {
        $codeShape = $slide2.Shapes.Item($i)
        break
    }
}

$tr = $codeShape.TextFrame.TextRange

$replacements = @(
    @{old="() = 1"; new="() -> 1"},
    @{old="()  = 1"; new="()  -> 1"},
    @{old="() = 1`t`t   "; new="() -> 1`t`t   "},
    @{old="()  = 1"; new="()  -> 1"},
    @{old="()  = 1"; new="() -> 1"},
    @{old="() = 2"; new="() -> 2"},
    @{old="() = 2"; new="() -> 2"},
    @{old="() = 3"; new="() -> 3"}
)

foreach ($r in $replacements) {
    $full = $tr.Text
    $idx = $full.IndexOf($r.old)
    if ($idx -ge 0) {
        $sub = $tr.Characters($idx + 1, $r.old.Length)
        $sub.Text = $r.new
    }
}
